# Regenerate Handback report timestamps for zh-cn and de-de sheets.
$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-24 15:30:40"
$wsZhCn.Range("H2").Value = "2016-03-24 15:31:19"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-24 15:30:45"
$wsDeDe.Range("H2").Value = "2016-03-24 15:31:28"
